$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3522
$ws.Range("I113").Value = 2544
$ws.Range("K113").Value = 2544
$ws.Range("M113").Value = 710

$ws.Range("H137").Value = 2191
$ws.Range("I137").Value = 2121.3333
$ws.Range("K137").Value = 6363.999899999999
$ws.Range("M137").Value = -3813.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 742.3333
$ws.Range("I2").Value = 582.5454999999999
$ws.Range("K2").Value = 582.5454999999999
$ws.Range("M2").Value = -469.5454999999999

$ws.Range("H32").Value = 6797.857
$ws.Range("I32").Value = 7189.1665
$ws.Range("K32").Value = 7189.1665
$ws.Range("M32").Value = -6902.1665

$ws.Range("H36").Value = 36666.332
$ws.Range("I36").Value = 36666.332
$ws.Range("K36").Value = 36666.332
$ws.Range("M36").Value = -36320.332

$ws.Range("H61").Value = 6489.5454
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 6489.5454
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 6489.5454
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -6913.5454

$ws.Range("H110").Value = 2847154.2
$ws.Range("I110").Value = 3364637.5
$ws.Range("K110").Value = 3364637.5
$ws.Range("M110").Value = -3362592.5

$ws.Range("H116").Value = 742.3333
$ws.Range("I116").Value = 582.5454999999999
$ws.Range("K116").Value = 582.5454999999999
$ws.Range("M116").Value = 1711.4545

$ws.Range("H136").Value = 6489.5454
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 6489.5454
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 19468.6362
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -24568.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 742.3333
$ws.Range("I3").Value = 582.5454999999999
$ws.Range("K3").Value = 582.5454999999999
$ws.Range("M3").Value = -468.5454999999999

$ws.Range("H54").Value = 8500
$ws.Range("I54").Value = 8500
$ws.Range("K54").Value = 8500
$ws.Range("M54").Value = -8016

$ws.Range("H76").Value = 106485.5
$ws.Range("J76").Value = 106485.5
$ws.Range("L76").Value = 106485.5
$ws.Range("N76").Value = -107115.5

$ws.Range("H79").Value = 106485.5
$ws.Range("J79").Value = 106485.5
$ws.Range("L79").Value = 106485.5
$ws.Range("N79").Value = -108669.5

$ws.Range("H88").Value = 42145
$ws.Range("J88").Value = 42145
$ws.Range("L88").Value = 42145
$ws.Range("N88").Value = -42957

$ws.Range("H91").Value = 42145
$ws.Range("J91").Value = 42145
$ws.Range("L91").Value = 42145
$ws.Range("N91").Value = -44953

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 733.3333
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = 150

$ws.Range("H31").Value = 1609.7059
$ws.Range("I31").Value = 1429.8
$ws.Range("J31").Value = 1684.6666
$ws.Range("K31").Value = 1429.8
$ws.Range("L31").Value = 1684.6666
$ws.Range("M31").Value = -1134.8
$ws.Range("N31").Value = -2274.6666

$ws.Range("H34").Value = 1609.7059
$ws.Range("I34").Value = 1429.8
$ws.Range("J34").Value = 1684.6666
$ws.Range("K34").Value = 1429.8
$ws.Range("L34").Value = 1684.6666
$ws.Range("M34").Value = -1227.8
$ws.Range("N34").Value = -2088.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9999
$ws.Range("I3").Value = 9999
$ws.Range("K3").Value = 29997
$ws.Range("M3").Value = -29885

$ws.Range("H37").Value = 63998.668
$ws.Range("J37").Value = 63998.668
$ws.Range("L37").Value = 191996.004
$ws.Range("N37").Value = -192220.004

$ws.Range("H68").Value = 13556.556
$ws.Range("J68").Value = 13556.556
$ws.Range("L68").Value = 40669.66800000001
$ws.Range("N68").Value = -42291.66800000001

$ws.Range("H71").Value = 13556.556
$ws.Range("J71").Value = 13556.556
$ws.Range("L71").Value = 122009.004
$ws.Range("N71").Value = -130121.004

$ws.Range("H98").Value = 217.71428
$ws.Range("J98").Value = 245.66667
$ws.Range("L98").Value = 737.00001
$ws.Range("N98").Value = -3733.00001

$ws.Range("H107").Value = 350
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 350
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1050
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -4890

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7992.2
$ws.Range("I70").Value = 7000.1
$ws.Range("J70").Value = 9976.4
$ws.Range("K70").Value = 7000.1
$ws.Range("L70").Value = 9976.4
$ws.Range("M70").Value = -6730.1
$ws.Range("N70").Value = -10516.4

$ws.Range("H73").Value = 7992.2
$ws.Range("I73").Value = 7000.1
$ws.Range("J73").Value = 9976.4
$ws.Range("K73").Value = 7000.1
$ws.Range("L73").Value = 9976.4
$ws.Range("M73").Value = -6064.1
$ws.Range("N73").Value = -11848.4

$ws.Range("H113").Value = 2862.8333
$ws.Range("I113").Value = 1340.2858
$ws.Range("K113").Value = 1340.2858
$ws.Range("M113").Value = 829.7141999999999

$ws.Range("H132").Value = 1841.2
$ws.Range("I132").Value = 1802
$ws.Range("K132").Value = 5406
$ws.Range("M132").Value = -2876

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 537.5
$ws.Range("I9").Value = 566.6667
$ws.Range("J9").Value = 450
$ws.Range("K9").Value = 566.6667
$ws.Range("L9").Value = 450
$ws.Range("M9").Value = -342.6667
$ws.Range("N9").Value = -898

$ws.Range("H22").Value = 4639.8335
$ws.Range("I22").Value = 4210.2
$ws.Range("K22").Value = 4210.2
$ws.Range("M22").Value = -3915.2

$ws.Range("H27").Value = 4639.8335
$ws.Range("I27").Value = 4210.2
$ws.Range("K27").Value = 4210.2
$ws.Range("M27").Value = -4103.2

$ws.Range("H122").Value = 6981.6855
$ws.Range("I122").Value = 6689.8696
$ws.Range("K122").Value = 20069.6088
$ws.Range("M122").Value = -17619.6088

$ws.Range("H132").Value = 2876.5
$ws.Range("I132").Value = 2407.5
$ws.Range("K132").Value = 7222.5
$ws.Range("M132").Value = -4692.5

$ws.Range("H136").Value = 4143.4287
$ws.Range("I136").Value = 3800.8
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 11402.4
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -8852.400000000001
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1076.125
$ws.Range("I113").Value = 1244.3334
$ws.Range("J113").Value = 571.5
$ws.Range("K113").Value = 3733.0002
$ws.Range("L113").Value = 1714.5
$ws.Range("M113").Value = -1563.0002
$ws.Range("N113").Value = -6054.5

$ws.Range("H132").Value = 7481.75
$ws.Range("I132").Value = 7198.3335
$ws.Range("K132").Value = 21595.0005
$ws.Range("M132").Value = -19065.0005
